# Refresh the cryptocurrency price/volume snapshot in column D (Price) and
# column E (Volume(1h)) for rows 2-51, per the "Updated cryptos list" run.
#
# Price values are stored as TEXT (not numbers) in the sheet -- e.g. prices
# like "30.628.02" use "." as a thousands separator, and plain decimals such
# as "1.013" must stay literal text too (Excel would otherwise silently parse
# them as numbers and round/reformat them). To force text we temporarily set
# NumberFormat to "@" before writing the value, then reset the cell's Style
# back to "Normal" so no stray formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value. Values that look like plain numbers
# ("1.013", "339.26", ...) are flagged so we can force-format them as text.
$updates = @(
    @{ Cell = "D2"; Value = "30.628.02"; ForceText = $false }
    @{ Cell = "D3"; Value = "2.112.23"; ForceText = $false }
    @{ Cell = "E3"; Value = "  +0.09%  "; ForceText = $false }
    @{ Cell = "D4"; Value = "1.013"; ForceText = $true }
    @{ Cell = "E4"; Value = "  +1.08%  "; ForceText = $false }
    @{ Cell = "D5"; Value = "339.26"; ForceText = $true }
    @{ Cell = "E6"; Value = "  +1.06%  "; ForceText = $false }
    @{ Cell = "D7"; Value = "0.5254"; ForceText = $true }
    @{ Cell = "E7"; Value = "  -0.10%  "; ForceText = $false }
    @{ Cell = "D8"; Value = "0.4497"; ForceText = $true }
    @{ Cell = "E8"; Value = "  +0.09%  "; ForceText = $false }
    @{ Cell = "D9"; Value = "53.70"; ForceText = $true }
    @{ Cell = "E9"; Value = "  +0.62%  "; ForceText = $false }
    @{ Cell = "D10"; Value = "0.09029"; ForceText = $true }
    @{ Cell = "E10"; Value = "  -0.02%  "; ForceText = $false }
    @{ Cell = "D11"; Value = "1.169"; ForceText = $true }
    @{ Cell = "E11"; Value = "  -0.81%  "; ForceText = $false }
    @{ Cell = "D12"; Value = "24.33"; ForceText = $true }
    @{ Cell = "E12"; Value = "  -0.60%  "; ForceText = $false }
    @{ Cell = "D13"; Value = "2.123.71"; ForceText = $false }
    @{ Cell = "E13"; Value = "  +1.11%  "; ForceText = $false }
    @{ Cell = "D14"; Value = "6.781"; ForceText = $true }
    @{ Cell = "E14"; Value = "  -0.11%  "; ForceText = $false }
    @{ Cell = "D15"; Value = "8.053"; ForceText = $true }
    @{ Cell = "E15"; Value = "  +2.95%  "; ForceText = $false }
    @{ Cell = "D16"; Value = "97.79"; ForceText = $true }
    @{ Cell = "E16"; Value = "  +1.05%  "; ForceText = $false }
    @{ Cell = "D17"; Value = "0.00001160"; ForceText = $true }
    @{ Cell = "E17"; Value = "  +2.59%  "; ForceText = $false }
    @{ Cell = "D18"; Value = "1.014"; ForceText = $true }
    @{ Cell = "E18"; Value = "  +1.01%  "; ForceText = $false }
    @{ Cell = "D19"; Value = "0.06694"; ForceText = $true }
    @{ Cell = "E19"; Value = "  +1.09%  "; ForceText = $false }
    @{ Cell = "D20"; Value = "19.34"; ForceText = $true }
    @{ Cell = "E20"; Value = "  -0.60%  "; ForceText = $false }
    @{ Cell = "D21"; Value = "1.012"; ForceText = $true }
    @{ Cell = "E21"; Value = "  +1.04%  "; ForceText = $false }
    @{ Cell = "D22"; Value = "6.379"; ForceText = $true }
    @{ Cell = "E22"; Value = "  +0.89%  "; ForceText = $false }
    @{ Cell = "D23"; Value = "30.724.52"; ForceText = $false }
    @{ Cell = "E23"; Value = "  +0.61%  "; ForceText = $false }
    @{ Cell = "E24"; Value = "  +3.11%  "; ForceText = $false }
    @{ Cell = "D25"; Value = "2.372"; ForceText = $true }
    @{ Cell = "E25"; Value = "  +0.68%  "; ForceText = $false }
    @{ Cell = "D26"; Value = "2.373.78"; ForceText = $false }
    @{ Cell = "E26"; Value = "  +1.14%  "; ForceText = $false }
    @{ Cell = "D27"; Value = "22.38"; ForceText = $true }
    @{ Cell = "E27"; Value = "  -0.10%  "; ForceText = $false }
    @{ Cell = "D28"; Value = "165.10"; ForceText = $true }
    @{ Cell = "E28"; Value = "  +0.94%  "; ForceText = $false }
    @{ Cell = "D29"; Value = "2.542"; ForceText = $true }
    @{ Cell = "E29"; Value = "  -1.86%  "; ForceText = $false }
    @{ Cell = "D30"; Value = "134.72"; ForceText = $true }
    @{ Cell = "E30"; Value = "  +1.37%  "; ForceText = $false }
    @{ Cell = "D31"; Value = "1.194"; ForceText = $true }
    @{ Cell = "E31"; Value = "  -0.57%  "; ForceText = $false }
    @{ Cell = "D32"; Value = "0.1074"; ForceText = $true }
    @{ Cell = "E32"; Value = "  -0.15%  "; ForceText = $false }
    @{ Cell = "D33"; Value = "6.363"; ForceText = $true }
    @{ Cell = "E33"; Value = "  +3.18%  "; ForceText = $false }
    @{ Cell = "D34"; Value = "1.625"; ForceText = $true }
    @{ Cell = "E34"; Value = "  -2.72%  "; ForceText = $false }
    @{ Cell = "D35"; Value = "3.942"; ForceText = $true }
    @{ Cell = "E35"; Value = "  +0.47%  "; ForceText = $false }
    @{ Cell = "D36"; Value = "10.31"; ForceText = $true }
    @{ Cell = "E36"; Value = "  -2.57%  "; ForceText = $false }
    @{ Cell = "D37"; Value = "5.884"; ForceText = $true }
    @{ Cell = "E37"; Value = "  +5.47%  "; ForceText = $false }
    @{ Cell = "D38"; Value = "0.02647"; ForceText = $true }
    @{ Cell = "E38"; Value = "  +2.37%  "; ForceText = $false }
    @{ Cell = "D39"; Value = "0.06821"; ForceText = $true }
    @{ Cell = "E39"; Value = "  -0.30%  "; ForceText = $false }
    @{ Cell = "D40"; Value = "0.2316"; ForceText = $true }
    @{ Cell = "E40"; Value = "  +0.32%  "; ForceText = $false }
    @{ Cell = "D41"; Value = "12.58"; ForceText = $true }
    @{ Cell = "E41"; Value = "  -1.69%  "; ForceText = $false }
    @{ Cell = "D42"; Value = "0.6866"; ForceText = $true }
    @{ Cell = "D43"; Value = "1.259"; ForceText = $true }
    @{ Cell = "E43"; Value = "  +0.93%  "; ForceText = $false }
    @{ Cell = "D44"; Value = "14.98"; ForceText = $true }
    @{ Cell = "E44"; Value = "  +6.51%  "; ForceText = $false }
    @{ Cell = "D45"; Value = "0.6418"; ForceText = $true }
    @{ Cell = "E45"; Value = "  +0.15%  "; ForceText = $false }
    @{ Cell = "D46"; Value = "2.312"; ForceText = $true }
    @{ Cell = "E46"; Value = "  -2.01%  "; ForceText = $false }
    @{ Cell = "D47"; Value = "0.00000000369"; ForceText = $true }
    @{ Cell = "E47"; Value = "  +11.03%  "; ForceText = $false }
    @{ Cell = "D48"; Value = "3.707"; ForceText = $true }
    @{ Cell = "E48"; Value = "  +1.16%  "; ForceText = $false }
    @{ Cell = "D49"; Value = "1.254"; ForceText = $true }
    @{ Cell = "D50"; Value = "82.88"; ForceText = $true }
    @{ Cell = "E50"; Value = "  -0.83%  "; ForceText = $false }
    @{ Cell = "D51"; Value = "0.07305"; ForceText = $true }
    @{ Cell = "E51"; Value = "  +3.07%  "; ForceText = $false }
)

foreach ($update in $updates) {
    $rng = $ws.Range($update.Cell)
    if ($update.ForceText) {
        # Looks numeric -- pin the format to Text first so Excel stores the
        # literal digits/trailing zeros instead of converting to a Double.
        $rng.NumberFormat = "@"
        $rng.Value = $update.Value
        $rng.Style = "Normal"
    } else {
        # Already unambiguous text (percent string, or multi-dot price) --
        # a plain assignment keeps it as text with no format juggling needed.
        $rng.Value = $update.Value
    }
}

